# Auto-generated Excel COM-interop script
# Applies numeric updates to the FFXIV Leve profit workbook sheets
# (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) as described by the source diff.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 31.5
$ws.Range("I5").Value = 44
$ws.Range("J5").Value = 19
$ws.Range("K5").Value = 44
$ws.Range("L5").Value = 19
$ws.Range("M5").Value = 71
$ws.Range("N5").Value = -249

$ws.Range("H13").Value = 5000
$ws.Range("I13").Value = 5000
$ws.Range("K13").Value = 5000
$ws.Range("M13").Value = -4831

$ws.Range("H17").Value = 2557.6
$ws.Range("J17").Value = 2557.6
$ws.Range("L17").Value = 7672.799999999999
$ws.Range("N17").Value = -8008.799999999999

$ws.Range("H76").Value = 3603.8333
$ws.Range("J76").Value = 3884.3333
$ws.Range("L76").Value = 3884.3333
$ws.Range("N76").Value = -4514.3333

$ws.Range("H79").Value = 3603.8333
$ws.Range("J79").Value = 3884.3333
$ws.Range("L79").Value = 3884.3333
$ws.Range("N79").Value = -6068.3333

$ws.Range("H141").Value = 2301.889
$ws.Range("I141").Value = 2005.9231
$ws.Range("K141").Value = 6017.7693
$ws.Range("M141").Value = -837.7692999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4407.1333
$ws.Range("I32").Value = 4407.1333
$ws.Range("K32").Value = 4407.1333
$ws.Range("M32").Value = -4120.1333

$ws.Range("H61").Value = 2601.4285
$ws.Range("I61").Value = 2601.4285
$ws.Range("K61").Value = 2601.4285
$ws.Range("M61").Value = -2389.4285

$ws.Range("H74").Value = 4602.2
$ws.Range("I74").Value = 3752.75
$ws.Range("K74").Value = 3752.75
$ws.Range("M74").Value = -2878.75

$ws.Range("H77").Value = 4602.2
$ws.Range("I77").Value = 3752.75
$ws.Range("K77").Value = 18763.75
$ws.Range("M77").Value = -14395.75

$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("N133").ClearContents()

$ws.Range("H136").Value = 2601.4285
$ws.Range("I136").Value = 2601.4285
$ws.Range("K136").Value = 7804.2855
$ws.Range("M136").Value = -5254.2855

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 392.2
$ws.Range("I22").Value = 392.2
$ws.Range("K22").Value = 392.2
$ws.Range("M22").Value = -219.2

$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()

$ws.Range("H105").Value = 2500
$ws.Range("I105").Value = 0
$ws.Range("J105").Value = 2500
$ws.Range("K105").Value = 0
$ws.Range("L105").Value = 2500
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -5994

$ws.Range("H134").Value = 4993.04
$ws.Range("I134").Value = 4993.04
$ws.Range("K134").Value = 14979.12
$ws.Range("M134").Value = -12444.12

$ws.Range("H137").Value = 41999.4
$ws.Range("I137").Value = 30000
$ws.Range("K137").Value = 30000
$ws.Range("M137").Value = -24900

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 65.75
$ws.Range("I7").Value = 65.75
$ws.Range("K7").Value = 65.75
$ws.Range("M7").Value = 47.25

$ws.Range("H13").Value = 9999
$ws.Range("I13").Value = 0
$ws.Range("K13").Value = 0
$ws.Range("M13").ClearContents()

$ws.Range("H22").Value = 99.333336
$ws.Range("I22").Value = 99.333336
$ws.Range("K22").Value = 99.333336
$ws.Range("M22").Value = 250.666664

$ws.Range("H31").Value = 4081.8235
$ws.Range("J31").Value = 6310.6
$ws.Range("L31").Value = 6310.6
$ws.Range("N31").Value = -6900.6

$ws.Range("H34").Value = 4081.8235
$ws.Range("J34").Value = 6310.6
$ws.Range("L34").Value = 6310.6
$ws.Range("N34").Value = -6714.6

$ws.Range("H58").Value = 2767.25
$ws.Range("I58").Value = 2767.25
$ws.Range("K58").Value = 2767.25
$ws.Range("M58").Value = -2564.25

$ws.Range("H92").Value = 29177.445
$ws.Range("J92").Value = 29137.125
$ws.Range("L92").Value = 29137.125
$ws.Range("N92").Value = -34129.125

$ws.Range("H132").Value = 2723.1667
$ws.Range("I132").Value = 1486
$ws.Range("J132").Value = 3341.75
$ws.Range("K132").Value = 4458
$ws.Range("L132").Value = 10025.25
$ws.Range("M132").Value = -1928
$ws.Range("N132").Value = -15085.25

$ws.Range("H136").Value = 2767.25
$ws.Range("I136").Value = 2767.25
$ws.Range("K136").Value = 8301.75
$ws.Range("M136").Value = -5751.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H117").Value = 906.3333
$ws.Range("I117").Value = 769
$ws.Range("J117").Value = 975
$ws.Range("K117").Value = 2307
$ws.Range("L117").Value = 2925
$ws.Range("M117").Value = 1135
$ws.Range("N117").Value = -9809

$ws.Range("H129").Value = 708.6667
$ws.Range("I129").Value = 708.6667
$ws.Range("K129").Value = 2126.0001
$ws.Range("M129").Value = 2873.9999

$ws.Range("H134").Value = 822
$ws.Range("I134").Value = 822
$ws.Range("K134").Value = 2466
$ws.Range("M134").Value = 2604

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").ClearContents()

$ws.Range("H99").Value = 5562.1665
$ws.Range("I99").Value = 5562.1665
$ws.Range("K99").Value = 5562.1665
$ws.Range("M99").Value = -3316.1665

$ws.Range("H132").Value = 2114.889
$ws.Range("I132").Value = 2114.889
$ws.Range("K132").Value = 6344.667
$ws.Range("M132").Value = -3814.667

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1450
$ws.Range("I22").Value = 2350
$ws.Range("K22").Value = 2350
$ws.Range("M22").Value = -2055

$ws.Range("H27").Value = 1450
$ws.Range("I27").Value = 2350
$ws.Range("K27").Value = 2350
$ws.Range("M27").Value = -2243

$ws.Range("H47").Value = 19514.75
$ws.Range("I47").Value = 5000
$ws.Range("K47").Value = 5000
$ws.Range("M47").Value = -4510

$ws.Range("H52").Value = 19514.75
$ws.Range("I52").Value = 5000
$ws.Range("K52").Value = 5000
$ws.Range("M52").Value = -4767

$ws.Range("H101").Value = 26181
$ws.Range("J101").Value = 26181
$ws.Range("L101").Value = 26181
$ws.Range("N101").Value = -32671

$ws.Range("H132").Value = 6207.5835
$ws.Range("I132").Value = 4998.5
$ws.Range("K132").Value = 14995.5
$ws.Range("M132").Value = -12465.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H127").Value = 50000.5
$ws.Range("J127").Value = 50000.5
$ws.Range("L127").Value = 50000.5
$ws.Range("N127").Value = -59920.5

$ws.Range("H132").Value = 1667.826
$ws.Range("I132").Value = 1478.9445
$ws.Range("K132").Value = 4436.833500000001
$ws.Range("M132").Value = -1906.833500000001

$ws.Range("H136").Value = 3404.1765
$ws.Range("I136").Value = 2626.5715
$ws.Range("K136").Value = 7879.7145
$ws.Range("M136").Value = -5329.7145

